$d = $word.ActiveDocument

# --- Recreate the style sheet so MSC_Join_A / MSC_Join_B / MSC_Join_C are
# --- interleaved with the existing MSC_Paragraph_A/B/C styles in the same
# --- order the target document uses. Styles.Add always appends at the end
# --- of the collection, so MSC_Paragraph_B and MSC_Paragraph_C are removed
# --- first (end-to-start, to avoid stale-index issues) and then re-added
# --- in the right spots, preserving their original definitions exactly.

$d.Styles.Item("MSC_Paragraph_C").Delete()
$d.Styles.Item("MSC_Paragraph_B").Delete()

$joinA = $d.Styles.Add("MSC_Join_A", 1)
$joinA.BaseStyle = "MSCJoin"

$paraB = $d.Styles.Add("MSC_Paragraph_B", 1)
$paraB.BaseStyle = "MSCParagraph"
$paraB.Font.NameAscii = "Noto Sans CJK TC"
$paraB.Font.Name = "Noto Sans CJK TC"
$paraB.Font.NameFarEast = "Noto Sans CJK TC"
$paraB.Font.NameBi = "Noto Sans CJK TC"

$joinB = $d.Styles.Add("MSC_Join_B", 1)
$joinB.BaseStyle = "MSCJoin"
$joinB.Font.NameAscii = "Noto Sans CJK TC"
$joinB.Font.Name = "Noto Sans CJK TC"
$joinB.Font.NameFarEast = "Noto Sans CJK TC"
$joinB.Font.NameBi = "Noto Sans CJK TC"

$paraC = $d.Styles.Add("MSC_Paragraph_C", 1)
$paraC.BaseStyle = "MSCParagraph"

$joinC = $d.Styles.Add("MSC_Join_C", 1)
$joinC.BaseStyle = "MSCJoin"

# --- Re-point the "join" paragraphs (the blank / "[...]" / blank triad that
# --- sits between two MSC_Paragraph_A or MSC_Paragraph_B verse blocks) from
# --- the shared MSCJoin style onto the new per-column MSC_Join_A / MSC_Join_B
# --- styles, matching each paragraph's column.

$joinAIndexes = @(9, 10, 11, 26, 27, 28)
foreach ($i in $joinAIndexes) {
    $d.Paragraphs.Item($i).Style = "MSC_Join_A"
}

$joinBIndexes = @(14, 15, 16, 31, 32, 33)
foreach ($i in $joinBIndexes) {
    $d.Paragraphs.Item($i).Style = "MSC_Join_B"
}
